$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 14520605
$ws.Range("I12").Value = 17424528
$ws.Range("K12").Value = 17424528
$ws.Range("M12").Value = -17424358
$ws.Range("H41").Value = 931.7692
$ws.Range("I41").Value = 453.5
$ws.Range("K41").Value = 453.5
$ws.Range("M41").Value = -13.5
$ws.Range("H76").Value = 4074.75
$ws.Range("I76").Value = 4199
$ws.Range("J76").Value = 4033.3333
$ws.Range("K76").Value = 4199
$ws.Range("L76").Value = 4033.3333
$ws.Range("M76").Value = -3884
$ws.Range("N76").Value = -4663.3333
$ws.Range("H79").Value = 4074.75
$ws.Range("I79").Value = 4199
$ws.Range("J79").Value = 4033.3333
$ws.Range("K79").Value = 4199
$ws.Range("L79").Value = 4033.3333
$ws.Range("M79").Value = -3107
$ws.Range("N79").Value = -6217.3333
$ws.Range("H111").Value = 1053.6
$ws.Range("I111").Value = 979.6667
$ws.Range("K111").Value = 2939.0001
$ws.Range("M111").Value = 127.9998999999998

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 797597.8
$ws.Range("I74").Value = 1272.5454
$ws.Range("K74").Value = 1272.5454
$ws.Range("M74").Value = -398.5454
$ws.Range("H77").Value = 797597.8
$ws.Range("I77").Value = 1272.5454
$ws.Range("K77").Value = 6362.727
$ws.Range("M77").Value = -1994.727
$ws.Range("H96").Value = 35608.8
$ws.Range("J96").Value = 35608.8
$ws.Range("L96").Value = 35608.8
$ws.Range("N96").Value = -41100.8
$ws.Range("H97").Value = 5315.96
$ws.Range("I97").Value = 6897.25
$ws.Range("J97").Value = 2504.7778
$ws.Range("K97").Value = 6897.25
$ws.Range("L97").Value = 2504.7778
$ws.Range("M97").Value = -6401.25
$ws.Range("N97").Value = -3496.7778

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10126.5
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 12835.333
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 12835.333
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -15081.333
$ws.Range("H89").Value = 10126.5
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 12835.333
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 64176.665
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -75408.66500000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 35000
$ws.Range("J28").Value = 35000
$ws.Range("L28").Value = 35000
$ws.Range("N28").Value = -35490
$ws.Range("H88").Value = 47671.25
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 47671.25
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 47671.25
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -48483.25
$ws.Range("H91").Value = 47671.25
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 47671.25
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 47671.25
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -50479.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 651.375
$ws.Range("I12").Value = 99.59999999999999
$ws.Range("K12").Value = 298.8
$ws.Range("M12").Value = -125.8
$ws.Range("H81").Value = 6788.75
$ws.Range("I81").Value = 5000
$ws.Range("J81").Value = 7385
$ws.Range("K81").Value = 15000
$ws.Range("L81").Value = 22155
$ws.Range("M81").Value = -13877
$ws.Range("N81").Value = -24401
$ws.Range("H84").Value = 6788.75
$ws.Range("I84").Value = 5000
$ws.Range("J84").Value = 7385
$ws.Range("K84").Value = 45000
$ws.Range("L84").Value = 66465
$ws.Range("M84").Value = -39384
$ws.Range("N84").Value = -77697
$ws.Range("H122").Value = 6960821
$ws.Range("I122").Value = 33333888
$ws.Range("J122").Value = 2165718.2
$ws.Range("K122").Value = 300004992
$ws.Range("L122").Value = 19491463.8
$ws.Range("M122").Value = -300002542
$ws.Range("N122").Value = -19496363.8
$ws.Range("H137").Value = 2714.1333
$ws.Range("I137").Value = 1589
$ws.Range("J137").Value = 4964.4
$ws.Range("K137").Value = 4767
$ws.Range("L137").Value = 14893.2
$ws.Range("M137").Value = 333
$ws.Range("N137").Value = -25093.2

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 4618.1
$ws.Range("I41").Value = 3464.5557
$ws.Range("K41").Value = 3464.5557
$ws.Range("M41").Value = -3109.5557
$ws.Range("H80").Value = 29427558
$ws.Range("I80").Value = 18095.715
$ws.Range("K80").Value = 18095.715
$ws.Range("M80").Value = -17097.715
$ws.Range("H83").Value = 29427558
$ws.Range("I83").Value = 18095.715
$ws.Range("K83").Value = 90478.575
$ws.Range("M83").Value = -85486.575
$ws.Range("H95").Value = 60000
$ws.Range("J95").Value = 60000
$ws.Range("L95").Value = 60000
$ws.Range("N95").Value = -65492
$ws.Range("H102").Value = 35716170
$ws.Range("I102").Value = 41668240
$ws.Range("K102").Value = 41668240
$ws.Range("M102").Value = -41666618
$ws.Range("H132").Value = 7069508
$ws.Range("I132").Value = 2291.3635
$ws.Range("K132").Value = 6874.0905
$ws.Range("M132").Value = -4344.0905

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1440.2106
$ws.Range("I16").Value = 1202.8235
$ws.Range("J16").Value = 3458
$ws.Range("K16").Value = 1202.8235
$ws.Range("L16").Value = 3458
$ws.Range("M16").Value = -1032.8235
$ws.Range("N16").Value = -3798
$ws.Range("H55").Value = 680.71875
$ws.Range("I55").Value = 321.66666
$ws.Range("K55").Value = 321.66666
$ws.Range("M55").Value = -148.66666
$ws.Range("H61").Value = 3155.9524
$ws.Range("I61").Value = 3010.3125
$ws.Range("K61").Value = 3010.3125
$ws.Range("M61").Value = -2808.3125
$ws.Range("H82").Value = 2079.125
$ws.Range("J82").Value = 2728
$ws.Range("L82").Value = 2728
$ws.Range("N82").Value = -3450
$ws.Range("H85").Value = 2079.125
$ws.Range("J85").Value = 2728
$ws.Range("L85").Value = 2728
$ws.Range("N85").Value = -5224
$ws.Range("H113").Value = 3155.9524
$ws.Range("I113").Value = 3010.3125
$ws.Range("K113").Value = 3010.3125
$ws.Range("M113").Value = -840.3125
$ws.Range("H122").Value = 4907.5
$ws.Range("I122").Value = 4868.143
$ws.Range("K122").Value = 14604.429
$ws.Range("M122").Value = -12154.429
$ws.Range("H132").Value = 3292.1177
$ws.Range("I132").Value = 3031
$ws.Range("K132").Value = 9093
$ws.Range("M132").Value = -6563

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 127841.875
$ws.Range("I81").Value = 3248
$ws.Range("J81").Value = 999999
$ws.Range("K81").Value = 6496
$ws.Range("L81").Value = 1999998
$ws.Range("M81").Value = -5435
$ws.Range("N81").Value = -2002120
$ws.Range("H84").Value = 127841.875
$ws.Range("I84").Value = 3248
$ws.Range("J84").Value = 999999
$ws.Range("K84").Value = 32480
$ws.Range("L84").Value = 9999990
$ws.Range("M84").Value = -27176
$ws.Range("N84").Value = -10010598
$ws.Range("H126").Value = 3578.2727
$ws.Range("I126").Value = 5396.75
$ws.Range("K126").Value = 16190.25
$ws.Range("M126").Value = -13720.25
$ws.Range("H132").Value = 2538.3333
$ws.Range("I132").Value = 2517.4666
$ws.Range("K132").Value = 7552.399800000001
$ws.Range("M132").Value = -5022.399800000001
